$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 9236.5
$ws.Range("I62").Value = 7500
$ws.Range("J62").Value = 9815.333000000001
$ws.Range("K62").Value = 7500
$ws.Range("L62").Value = 9815.333000000001
$ws.Range("M62").Value = -6876
$ws.Range("N62").Value = -11063.333
$ws.Range("H65").Value = 9236.5
$ws.Range("I65").Value = 7500
$ws.Range("J65").Value = 9815.333000000001
$ws.Range("K65").Value = 37500
$ws.Range("L65").Value = 49076.665
$ws.Range("M65").Value = -34380
$ws.Range("N65").Value = -55316.665
$ws.Range("H112").Value = 1057.5588
$ws.Range("J112").Value = 1108.5667
$ws.Range("L112").Value = 3325.7001
$ws.Range("N112").Value = -5541.7001
$ws.Range("H118").Value = 392.5
$ws.Range("I118").Value = 392.5
$ws.Range("K118").Value = 1177.5
$ws.Range("M118").Value = 479.5
$ws.Range("H132").Value = 2819.3333
$ws.Range("I132").Value = 3004.3428
$ws.Range("K132").Value = 9013.028399999999
$ws.Range("M132").Value = -6483.028399999999
$ws.Range("H135").Value = 29415010
$ws.Range("I135").Value = 1295.4286
$ws.Range("J135").Value = 166679000
$ws.Range("K135").Value = 11658.8574
$ws.Range("L135").Value = 1500111000
$ws.Range("M135").Value = -9123.857399999999
$ws.Range("N135").Value = -1500116070

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 685.1539
$ws.Range("I2").Value = 752.7368
$ws.Range("K2").Value = 752.7368
$ws.Range("M2").Value = -639.7368
$ws.Range("H32").Value = 6592.2974
$ws.Range("I32").Value = 4518.807
$ws.Range("J32").Value = 13544.588
$ws.Range("K32").Value = 4518.807
$ws.Range("L32").Value = 13544.588
$ws.Range("M32").Value = -4231.807
$ws.Range("N32").Value = -14118.588
$ws.Range("H61").Value = 1992.0358
$ws.Range("I61").Value = 1355
$ws.Range("J61").Value = 3584.625
$ws.Range("K61").Value = 1355
$ws.Range("L61").Value = 3584.625
$ws.Range("M61").Value = -1143
$ws.Range("N61").Value = -4008.625
$ws.Range("H116").Value = 685.1539
$ws.Range("I116").Value = 752.7368
$ws.Range("K116").Value = 752.7368
$ws.Range("M116").Value = 1541.2632
$ws.Range("H122").Value = 3341.75
$ws.Range("I122").Value = 2390.7144
$ws.Range("K122").Value = 7172.1432
$ws.Range("M122").Value = -4722.1432
$ws.Range("H132").Value = 14369.743
$ws.Range("I132").Value = 1313.1875
$ws.Range("K132").Value = 3939.5625
$ws.Range("M132").Value = -1409.5625
$ws.Range("H136").Value = 1992.0358
$ws.Range("I136").Value = 1355
$ws.Range("J136").Value = 3584.625
$ws.Range("K136").Value = 4065
$ws.Range("L136").Value = 10753.875
$ws.Range("M136").Value = -1515
$ws.Range("N136").Value = -15853.875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 685.1539
$ws.Range("I3").Value = 752.7368
$ws.Range("K3").Value = 752.7368
$ws.Range("M3").Value = -638.7368
$ws.Range("H105").Value = 1810
$ws.Range("I105").Value = 1560.909
$ws.Range("J105").Value = 1886.1111
$ws.Range("K105").Value = 1560.909
$ws.Range("L105").Value = 1886.1111
$ws.Range("M105").Value = 186.0909999999999
$ws.Range("N105").Value = -5380.1111
$ws.Range("H107").Value = 849.5833
$ws.Range("I107").Value = 820.2222
$ws.Range("J107").Value = 937.6667
$ws.Range("K107").Value = 820.2222
$ws.Range("L107").Value = 937.6667
$ws.Range("M107").Value = 1099.7778
$ws.Range("N107").Value = -4777.6667
$ws.Range("H134").Value = 3198.1082
$ws.Range("I134").Value = 3242.1177
$ws.Range("J134").Value = 2699.3333
$ws.Range("K134").Value = 9726.3531
$ws.Range("L134").Value = 8097.999899999999
$ws.Range("M134").Value = -7191.3531
$ws.Range("N134").Value = -13167.9999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4763.174
$ws.Range("I31").Value = 5888.8
$ws.Range("J31").Value = 4450.5
$ws.Range("K31").Value = 5888.8
$ws.Range("L31").Value = 4450.5
$ws.Range("M31").Value = -5593.8
$ws.Range("N31").Value = -5040.5
$ws.Range("H34").Value = 4763.174
$ws.Range("I34").Value = 5888.8
$ws.Range("J34").Value = 4450.5
$ws.Range("K34").Value = 5888.8
$ws.Range("L34").Value = 4450.5
$ws.Range("M34").Value = -5686.8
$ws.Range("N34").Value = -4854.5
$ws.Range("H99").Value = 20003208
$ws.Range("I99").Value = 2734.2222
$ws.Range("J99").Value = 71433000
$ws.Range("K99").Value = 2734.2222
$ws.Range("L99").Value = 71433000
$ws.Range("M99").Value = -1236.2222
$ws.Range("N99").Value = -71435996
$ws.Range("H126").Value = 20003208
$ws.Range("I126").Value = 2734.2222
$ws.Range("J126").Value = 71433000
$ws.Range("K126").Value = 8202.6666
$ws.Range("L126").Value = 214299000
$ws.Range("M126").Value = -5732.6666
$ws.Range("N126").Value = -214303940
$ws.Range("H132").Value = 3471.9524
$ws.Range("I132").Value = 2422.3333
$ws.Range("K132").Value = 7266.999899999999
$ws.Range("M132").Value = -4736.999899999999
$ws.Range("H134").Value = 1714.8667
$ws.Range("I134").Value = 1643.5834
$ws.Range("K134").Value = 4930.7502
$ws.Range("M134").Value = -2395.7502

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H136").Value = 2475.389
$ws.Range("I136").Value = 981.0909
$ws.Range("K136").Value = 2943.2727
$ws.Range("M136").Value = 2156.7273

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 41.954544
$ws.Range("I2").Value = 40.6
$ws.Range("J2").Value = 44.857143
$ws.Range("K2").Value = 40.6
$ws.Range("L2").Value = 44.857143
$ws.Range("M2").Value = 72.40000000000001
$ws.Range("N2").Value = -270.857143
$ws.Range("H126").Value = 3507.2273
$ws.Range("I126").Value = 2372.4375
$ws.Range("J126").Value = 6533.3335
$ws.Range("K126").Value = 7117.3125
$ws.Range("L126").Value = 19600.0005
$ws.Range("M126").Value = -4647.3125
$ws.Range("N126").Value = -24540.0005
$ws.Range("H132").Value = 27539.35
$ws.Range("I132").Value = 2071.5
$ws.Range("K132").Value = 6214.5
$ws.Range("M132").Value = -3684.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1951.8096
$ws.Range("I22").Value = 2240
$ws.Range("J22").Value = 727
$ws.Range("K22").Value = 2240
$ws.Range("L22").Value = 727
$ws.Range("M22").Value = -1945
$ws.Range("N22").Value = -1317
$ws.Range("H27").Value = 1951.8096
$ws.Range("I27").Value = 2240
$ws.Range("J27").Value = 727
$ws.Range("K27").Value = 2240
$ws.Range("L27").Value = 727
$ws.Range("M27").Value = -2133
$ws.Range("N27").Value = -941
$ws.Range("H122").Value = 983613.9
$ws.Range("I122").Value = 1403213.8
$ws.Range("J122").Value = 4547.5
$ws.Range("K122").Value = 4209641.4
$ws.Range("L122").Value = 13642.5
$ws.Range("M122").Value = -4207191.4
$ws.Range("N122").Value = -18542.5
$ws.Range("H132").Value = 604893.4
$ws.Range("I132").Value = 1507248.9
$ws.Range("J132").Value = 3323.0833
$ws.Range("K132").Value = 4521746.699999999
$ws.Range("L132").Value = 9969.249899999999
$ws.Range("M132").Value = -4519216.699999999
$ws.Range("N132").Value = -15029.2499
$ws.Range("H136").Value = 1788.2778
$ws.Range("I136").Value = 1729.2858
$ws.Range("K136").Value = 5187.857400000001
$ws.Range("M136").Value = -2637.857400000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 55024470
$ws.Range("I107").Value = 90909260
$ws.Range("J107").Value = 5682878
$ws.Range("K107").Value = 272727780
$ws.Range("L107").Value = 17048634
$ws.Range("M107").Value = -272725860
$ws.Range("N107").Value = -17052474
$ws.Range("H132").Value = 950.4286
$ws.Range("I132").Value = 603.3823
$ws.Range("J132").Value = 2425.375
$ws.Range("K132").Value = 1810.1469
$ws.Range("L132").Value = 7276.125
$ws.Range("M132").Value = 719.8531
$ws.Range("N132").Value = -12336.125
$ws.Range("H136").Value = 24579084
$ws.Range("I136").Value = 31281740
$ws.Range("J136").Value = 2677.7778
$ws.Range("K136").Value = 93845220
$ws.Range("L136").Value = 8033.3334
$ws.Range("M136").Value = -93842670
$ws.Range("N136").Value = -13133.3334
